$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old expense-tracker content that lived in A1:B5
$ws.Range("A1:B5").Clear() | Out-Null

# Export data into xlsx: a small product/size table
$data = @(
    @($false, "Prod1", "Size", "L"),
    @($false, $false,  $false, $false),
    @($false, $false,  $false, $false),
    @($false, $false,  $false, $false),
    @($false, $false,  $false, $false),
    @($false, $false,  $false, $false),
    @($false, $false,  $false, $false),
    @($false, "Prod2", "Size", "L"),
    @($false, "Prod2", "Size", "L"),
    @($false, "Prod1", "Size", $false),
    @("All",  "Prod1", "Size", "L")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
